$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

    $ws.Range("G2").Value = 14.062694
    $ws.Range("H2").Value = 42.188082
    $ws.Range("I2").Value = 0.2613715759657495
    $ws.Range("J2").Value = 0.2613715759657494
    $ws.Range("K2").Value = 2
    $ws.Range("L2").Value = 0.6666666666666666
    $ws.Range("M2").Value = 0.01650666666666667
    $ws.Range("N2").Value = 0.04952
    $ws.Range("O2").Value = 0.795859985214233
    $ws.Range("P2").Value = 0.795859985214233
    $ws.Range("Q2").Value = 0.2321282022933333
    $ws.Range("R2").Value = 2.08915382064
    $ws.Range("S2").Value = 0.2080151785835222
    $ws.Range("T2").Value = 0.2080151785835221
    $ws.Range("G3").Value = 14.062694
    $ws.Range("H3").Value = 42.188082
    $ws.Range("I3").Value = 0.2613715759657495
    $ws.Range("J3").Value = 0.2613715759657494
    $ws.Range("O3").Value = 0.2041400147857671
    $ws.Range("P3").Value = 0.2041400147857671
    $ws.Range("Q3").Value = 0.059541446396
    $ws.Range("R3").Value = 0.535873017564
    $ws.Range("S3").Value = 0.05335639738222735
    $ws.Range("T3").Value = 0.05335639738222734
    $ws.Range("I4").Value = 0.1749200409165788
    $ws.Range("J4").Value = 0.1749200409165788
    $ws.Range("K4").Value = 2
    $ws.Range("L4").Value = 0.6666666666666666
    $ws.Range("M4").Value = 0.01650666666666667
    $ws.Range("N4").Value = 0.04952
    $ws.Range("O4").Value = 0.795859985214233
    $ws.Range("P4").Value = 0.795859985214233
    $ws.Range("Q4").Value = 0.1553492360177778
    $ws.Range("R4").Value = 1.39814312416
    $ws.Range("S4").Value = 0.1392118611775414
    $ws.Range("T4").Value = 0.1392118611775414
    $ws.Range("I5").Value = 0.1749200409165788
    $ws.Range("J5").Value = 0.1749200409165788
    $ws.Range("O5").Value = 0.2041400147857671
    $ws.Range("P5").Value = 0.2041400147857671
    $ws.Range("S5").Value = 0.03570817973903738
    $ws.Range("T5").Value = 0.03570817973903738
    $ws.Range("H6").Value = 90.988377
    $ws.Range("I6").Value = 0.5637083831176718
    $ws.Range("J6").Value = 0.5637083831176717
    $ws.Range("K6").Value = 2
    $ws.Range("L6").Value = 0.6666666666666666
    $ws.Range("M6").Value = 0.01650666666666667
    $ws.Range("N6").Value = 0.04952
    $ws.Range("O6").Value = 0.795859985214233
    $ws.Range("P6").Value = 0.795859985214233
    $ws.Range("Q6").Value = 0.5006382698933333
    $ws.Range("R6").Value = 4.50574442904
    $ws.Range("S6").Value = 0.4486329454531695
    $ws.Range("T6").Value = 0.4486329454531693
    $ws.Range("H7").Value = 90.988377
    $ws.Range("I7").Value = 0.5637083831176718
    $ws.Range("J7").Value = 0.5637083831176717
    $ws.Range("O7").Value = 0.2041400147857671
    $ws.Range("P7").Value = 0.2041400147857671
    $ws.Range("S7").Value = 0.1150754376645024
    $ws.Range("T7").Value = 0.1150754376645024
